$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.434.92"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").Value = "1.827.47"
$ws.Range("E3").Value = "  +2.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.81"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5350"
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4027"
$ws.Range("E8").Value = "  +6.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07689"
$ws.Range("E9").Value = "  +3.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.87"
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.113"
$ws.Range("E11").Value = "  +2.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.324"
$ws.Range("E12").Value = "  +4.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.627"
$ws.Range("E13").Value = "  +5.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.98"
$ws.Range("E14").Value = "  +2.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.000"
$ws.Range("D16").Value = "1.829.53"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.73"
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001078"
$ws.Range("E18").Value = "  +2.45%  "
$ws.Range("E19").Value = "  +2.31%  "
$ws.Range("E20").Value = "  +2.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.084"
$ws.Range("E22").Value = "  +3.52%  "
$ws.Range("D23").Value = "28.448.68"
$ws.Range("E23").Value = "  +1.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.17"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.229"
$ws.Range("E25").Value = "  +6.81%  "
$ws.Range("E26").Value = "  +8.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.74"
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.73"
$ws.Range("E28").Value = "  +2.54%  "
$ws.Range("D29").Value = "2.040.62"
$ws.Range("E29").Value = "  +2.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.19"
$ws.Range("E30").Value = "  +3.66%  "
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1109"
$ws.Range("E32").Value = "  +5.24%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.677"
$ws.Range("E33").Value = "  +2.85%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07498"
$ws.Range("E34").Value = "  +16.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.642"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2241"
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02343"
$ws.Range("E37").Value = "  +2.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.231"
$ws.Range("E38").Value = "  +4.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.891"
$ws.Range("E39").Value = "  +5.67%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6275"
$ws.Range("E40").Value = "  +2.34%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.32"
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.177"
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9999"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.393"
$ws.Range("E44").Value = "  -3.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.57"
$ws.Range("E45").Value = "  +2.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.700"
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5855"
$ws.Range("E47").Value = "  +2.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.92"
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.006"
$ws.Range("E49").Value = "  +4.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.206"
$ws.Range("E50").Value = "  +1.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06896"
$ws.Range("E51").Value = "  +1.57%  "
